$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 88, shifting the existing rows 88-122 down to 89-123.
$ws.Rows.Item(88).Insert()

# Populate the newly inserted row 88 with this week's record (same dims/variety
# as the prior top record, new date + prices).
$ws.Cells.Item(88, 1).Value = 7
$ws.Cells.Item(88, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(88, 3).Value = "Ñuble"
$ws.Cells.Item(88, 4).Value = 44460
$ws.Cells.Item(88, 5).Value = 16
$ws.Cells.Item(88, 6).Value = 100112003
$ws.Cells.Item(88, 7).Value = "Ajo"
$ws.Cells.Item(88, 8).Value = "Chino"
$ws.Cells.Item(88, 9).Value = "Primera"
$ws.Cells.Item(88, 10).Value = 120
$ws.Cells.Item(88, 11).Value = 16000
$ws.Cells.Item(88, 12).Value = 17000
$ws.Cells.Item(88, 13).Value = 16500
$ws.Cells.Item(88, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(88, 15).Value = "China"
$ws.Cells.Item(88, 16).Value = 1650
$ws.Cells.Item(88, 17).Value = 10
$ws.Cells.Item(88, 18).Value = "Hortaliza"
